$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 25 de Junio de 2020 a las 10:24"
$ws.Cells.Item(10, 3).Value = 4423
$ws.Cells.Item(10, 4).Value = 10332
$ws.Cells.Item(10, 5).Value = 7124
$ws.Cells.Item(12, 3).Value = 3368
$ws.Cells.Item(12, 4).Value = 6464
$ws.Cells.Item(12, 5).Value = 2030
$ws.Cells.Item(14, 3).Value = 2767
$ws.Cells.Item(14, 4).Value = 5609
$ws.Cells.Item(14, 5).Value = 4907
$ws.Cells.Item(15, 3).Value = 1812
$ws.Cells.Item(15, 4).Value = 5287
$ws.Cells.Item(15, 5).Value = 2820
$ws.Cells.Item(17, 3).Value = 4423
$ws.Cells.Item(17, 4).Value = 4868
$ws.Cells.Item(17, 5).Value = 7124
$ws.Cells.Item(18, 3).Value = 2454
$ws.Cells.Item(18, 4).Value = 4393
$ws.Cells.Item(18, 5).Value = 1577
$ws.Cells.Item(19, 3).Value = 2613
$ws.Cells.Item(19, 4).Value = 4152
$ws.Cells.Item(19, 5).Value = 1181
$ws.Cells.Item(21, 3).Value = 1124
$ws.Cells.Item(21, 4).Value = 3872
$ws.Cells.Item(21, 5).Value = 1992
$ws.Cells.Item(22, 3).Value = 1938
$ws.Cells.Item(22, 4).Value = 3794
$ws.Cells.Item(22, 5).Value = 3637
$ws.Cells.Item(23, 3).Value = 1872
$ws.Cells.Item(23, 4).Value = 3775
$ws.Cells.Item(23, 5).Value = 1389
$ws.Cells.Item(24, 3).Value = 1533
$ws.Cells.Item(24, 4).Value = 3569
$ws.Cells.Item(24, 5).Value = 1629
$ws.Cells.Item(25, 3).Value = 2344
$ws.Cells.Item(25, 4).Value = 3413
$ws.Cells.Item(25, 5).Value = 868
$ws.Cells.Item(26, 3).Value = 4423
$ws.Cells.Item(26, 4).Value = 3116
$ws.Cells.Item(26, 5).Value = 7124
$ws.Cells.Item(28, 3).Value = 185
$ws.Cells.Item(28, 4).Value = 2758
$ws.Cells.Item(28, 5).Value = 2295
$ws.Cells.Item(29, 3).Value = 1640
$ws.Cells.Item(29, 4).Value = 2746
$ws.Cells.Item(29, 5).Value = 901
$ws.Cells.Item(31, 3).Value = 308
$ws.Cells.Item(31, 4).Value = 2423
$ws.Cells.Item(31, 5).Value = 1840
$ws.Cells.Item(32, 4).Value = 2413
$ws.Cells.Item(32, 5).Value = 2520
$ws.Cells.Item(33, 1).Value = "Gran Canaria"
$ws.Cells.Item(33, 2).Value = 2347
$ws.Cells.Item(33, 3).Value = 659
$ws.Cells.Item(33, 4).Value = 2347
$ws.Cells.Item(33, 5).Value = 1537
$ws.Cells.Item(34, 1).Value = "Soria"
$ws.Cells.Item(34, 2).Value = 2290
$ws.Cells.Item(34, 3).Value = 1774
$ws.Cells.Item(34, 4).Value = 2290
$ws.Cells.Item(34, 5).Value = 397
$ws.Cells.Item(35, 1).Value = "Tenerife"
$ws.Cells.Item(35, 2).Value = 2280
$ws.Cells.Item(35, 3).Value = 623
$ws.Cells.Item(35, 4).Value = 2280
$ws.Cells.Item(35, 5).Value = 1506
$ws.Cells.Item(36, 1).Value = "Cantabria"
$ws.Cells.Item(36, 2).Value = 2246
$ws.Cells.Item(36, 3).Value = 1981
$ws.Cells.Item(36, 4).Value = 62
$ws.Cells.Item(36, 5).Value = 203
$ws.Cells.Item(37, 1).Value = "Caceres"
$ws.Cells.Item(37, 2).Value = 1973
$ws.Cells.Item(37, 3).Value = 66
$ws.Cells.Item(37, 4).Value = 1973
$ws.Cells.Item(37, 5).Value = 1505
$ws.Cells.Item(38, 1).Value = "A Coruña"
$ws.Cells.Item(38, 2).Value = 1969
$ws.Cells.Item(38, 3).Value = 333
$ws.Cells.Item(38, 4).Value = 1788
$ws.Cells.Item(38, 5).Value = 67
$ws.Cells.Item(39, 1).Value = "Avila"
$ws.Cells.Item(39, 2).Value = 1935
$ws.Cells.Item(39, 3).Value = 1179
$ws.Cells.Item(39, 4).Value = 1935
$ws.Cells.Item(39, 5).Value = 623
$ws.Cells.Item(40, 1).Value = "Murcia"
$ws.Cells.Item(40, 2).Value = 1587
$ws.Cells.Item(40, 3).Value = 2180
$ws.Cells.Item(40, 4).Value = 0
$ws.Cells.Item(40, 5).Value = 148
$ws.Cells.Item(41, 1).Value = "Pontevedra"
$ws.Cells.Item(41, 2).Value = 1536
$ws.Cells.Item(41, 3).Value = 333
$ws.Cells.Item(41, 4).Value = 1411
$ws.Cells.Item(41, 5).Value = 30
$ws.Cells.Item(42, 1).Value = "Castello/Castellon"
$ws.Cells.Item(42, 2).Value = 1486
$ws.Cells.Item(42, 3).Value = 699
$ws.Cells.Item(42, 4).Value = 1486
$ws.Cells.Item(42, 5).Value = 1363
$ws.Cells.Item(43, 1).Value = "Jaen"
$ws.Cells.Item(43, 2).Value = 1387
$ws.Cells.Item(43, 3).Value = 41
$ws.Cells.Item(43, 4).Value = 1387
$ws.Cells.Item(43, 5).Value = 1171
$ws.Cells.Item(44, 1).Value = "Cordoba"
$ws.Cells.Item(44, 2).Value = 1331
$ws.Cells.Item(44, 4).Value = 1331
$ws.Cells.Item(44, 5).Value = 1350
$ws.Cells.Item(45, 1).Value = "Guadalajara"
$ws.Cells.Item(45, 2).Value = 1266
$ws.Cells.Item(45, 3).Value = 644
$ws.Cells.Item(45, 4).Value = 1266
$ws.Cells.Item(45, 5).Value = 371
$ws.Cells.Item(46, 1).Value = "Cuenca"
$ws.Cells.Item(46, 2).Value = 1241
$ws.Cells.Item(46, 3).Value = 339
$ws.Cells.Item(46, 4).Value = 1241
$ws.Cells.Item(46, 5).Value = 596
$ws.Cells.Item(47, 1).Value = "Cadiz"
$ws.Cells.Item(47, 2).Value = 1240
$ws.Cells.Item(47, 3).Value = 535
$ws.Cells.Item(47, 4).Value = 1240
$ws.Cells.Item(47, 5).Value = 560
$ws.Cells.Item(48, 1).Value = "Palencia"
$ws.Cells.Item(48, 2).Value = 1205
$ws.Cells.Item(48, 3).Value = 789
$ws.Cells.Item(48, 4).Value = 1205
$ws.Cells.Item(48, 5).Value = 333
$ws.Cells.Item(49, 1).Value = "Huesca"
$ws.Cells.Item(49, 2).Value = 1115
$ws.Cells.Item(49, 3).Value = 544
$ws.Cells.Item(49, 4).Value = 1115
$ws.Cells.Item(49, 5).Value = 472
$ws.Cells.Item(50, 1).Value = "Zamora"
$ws.Cells.Item(50, 2).Value = 993
$ws.Cells.Item(50, 3).Value = 586
$ws.Cells.Item(50, 4).Value = 993
$ws.Cells.Item(50, 5).Value = 322
$ws.Cells.Item(51, 1).Value = "Badajoz"
$ws.Cells.Item(51, 2).Value = 962
$ws.Cells.Item(51, 3).Value = 0
$ws.Cells.Item(51, 4).Value = 962
$ws.Cells.Item(51, 5).Value = 1082
$ws.Cells.Item(52, 1).Value = "Ourense"
$ws.Cells.Item(52, 2).Value = 751
$ws.Cells.Item(52, 3).Value = 333
$ws.Cells.Item(52, 4).Value = 660
$ws.Cells.Item(52, 5).Value = 22
$ws.Cells.Item(53, 1).Value = "Teruel"
$ws.Cells.Item(53, 2).Value = 664
$ws.Cells.Item(53, 3).Value = 203
$ws.Cells.Item(53, 4).Value = 664
$ws.Cells.Item(53, 5).Value = 378
$ws.Cells.Item(54, 1).Value = "Lugo"
$ws.Cells.Item(54, 2).Value = 586
$ws.Cells.Item(54, 3).Value = 333
$ws.Cells.Item(54, 4).Value = 520
$ws.Cells.Item(54, 5).Value = 11
$ws.Cells.Item(55, 4).Value = 498
$ws.Cells.Item(55, 5).Value = 484
$ws.Cells.Item(56, 4).Value = 400
$ws.Cells.Item(56, 5).Value = 391
$ws.Cells.Item(60, 3).Value = 22
$ws.Cells.Item(60, 4).Value = 95
$ws.Cells.Item(60, 5).Value = 68
$ws.Cells.Item(61, 3).Value = 7
$ws.Cells.Item(61, 4).Value = 84
$ws.Cells.Item(61, 5).Value = 71
$ws.Cells.Item(63, 4).Value = 23
$ws.Cells.Item(63, 5).Value = 42
$ws.Cells.Item(66, 4).Value = 8
$ws.Cells.Item(68, 4).Value = 3

Write-Output "Done applying updates"
